$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 279, shifting existing row 279 (and all
# rows below it) down by one row.
$ws.Rows("279:279").Insert()

# Populate the newly inserted row 279 with the new record's data.
$ws.Range("A279").Value = 9
$ws.Range("B279").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C279").Value = 'Metropolitana'
$ws.Range("D279").Value = 44900
$ws.Range("D279").NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Range("E279").Value = 13
$ws.Range("F279").Value = 100112043
$ws.Range("G279").Value = 'Pepino ensalada'
$ws.Range("H279").Value = 'Sin especificar'
$ws.Range("I279").Value = 'Primera'
$ws.Range("J279").Value = 70
$ws.Range("K279").Value = 15000
$ws.Range("L279").Value = 16000
$ws.Range("M279").Value = 15500
$ws.Range("N279").Value = '$/caja 60 unidades'
$ws.Range("O279").Value = 'Región de Arica y Parinacota'
$ws.Range("P279").Value = 258
$ws.Range("Q279").Value = 60
$ws.Range("R279").Value = 'Hortaliza'
